$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 872.6
$ws.Range("I33").Value = 1062.6666
$ws.Range("J33").Value = 587.5
$ws.Range("K33").Value = 1062.6666
$ws.Range("L33").Value = 587.5
$ws.Range("M33").Value = -833.6666
$ws.Range("N33").Value = -1045.5
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 12000
$ws.Range("M46").Value = -11881
$ws.Range("H60").Value = 4000
$ws.Range("I60").Value = 4000
$ws.Range("K60").Value = 12000
$ws.Range("M60").Value = -11516
$ws.Range("H62").Value = 3981.3794
$ws.Range("I62").Value = 3445.1333
$ws.Range("J62").Value = 4555.9287
$ws.Range("K62").Value = 3445.1333
$ws.Range("L62").Value = 4555.9287
$ws.Range("M62").Value = -2821.1333
$ws.Range("N62").Value = -5803.9287
$ws.Range("H64").Value = 4677.8125
$ws.Range("I64").Value = 4484.5
$ws.Range("K64").Value = 4484.5
$ws.Range("M64").Value = -4236.5
$ws.Range("H65").Value = 3981.3794
$ws.Range("I65").Value = 3445.1333
$ws.Range("J65").Value = 4555.9287
$ws.Range("K65").Value = 17225.6665
$ws.Range("L65").Value = 22779.6435
$ws.Range("M65").Value = -14105.6665
$ws.Range("N65").Value = -29019.6435
$ws.Range("H67").Value = 4677.8125
$ws.Range("I67").Value = 4484.5
$ws.Range("K67").Value = 4484.5
$ws.Range("M67").Value = -3626.5
$ws.Range("H69").Value = 14495
$ws.Range("I69").Value = 3976.3333
$ws.Range("J69").Value = 16598.732
$ws.Range("K69").Value = 11928.9999
$ws.Range("L69").Value = 49796.196
$ws.Range("M69").Value = -11054.9999
$ws.Range("N69").Value = -51544.196
$ws.Range("H72").Value = 14495
$ws.Range("I72").Value = 3976.3333
$ws.Range("J72").Value = 16598.732
$ws.Range("K72").Value = 35786.9997
$ws.Range("L72").Value = 149388.588
$ws.Range("M72").Value = -31418.9997
$ws.Range("N72").Value = -158124.588
$ws.Range("H76").Value = 8156.5
$ws.Range("I76").Value = 7422.375
$ws.Range("K76").Value = 7422.375
$ws.Range("M76").Value = -7107.375
$ws.Range("H79").Value = 8156.5
$ws.Range("I79").Value = 7422.375
$ws.Range("K79").Value = 7422.375
$ws.Range("M79").Value = -6330.375
$ws.Range("H80").Value = 964.2143
$ws.Range("I80").Value = 891.5
$ws.Range("J80").Value = 1018.75
$ws.Range("K80").Value = 2674.5
$ws.Range("L80").Value = 3056.25
$ws.Range("M80").Value = -1676.5
$ws.Range("N80").Value = -5052.25
$ws.Range("H83").Value = 964.2143
$ws.Range("I83").Value = 891.5
$ws.Range("J83").Value = 1018.75
$ws.Range("K83").Value = 8023.5
$ws.Range("L83").Value = 9168.75
$ws.Range("M83").Value = -3031.5
$ws.Range("N83").Value = -19152.75
$ws.Range("H138").Value = 3139.2322
$ws.Range("I138").Value = 2573.5186
$ws.Range("K138").Value = 7720.5558
$ws.Range("M138").Value = -2580.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2840.2563
$ws.Range("I2").Value = 1479.1364
$ws.Range("K2").Value = 1479.1364
$ws.Range("M2").Value = -1366.1364
$ws.Range("H61").Value = 50457950
$ws.Range("I61").Value = 50457950
$ws.Range("K61").Value = 50457950
$ws.Range("M61").Value = -50457738
$ws.Range("H86").Value = 49999
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 49999
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 49999
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -52371
$ws.Range("H89").Value = 49999
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 49999
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 149997
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -161853
$ws.Range("H104").Value = 20142.666
$ws.Range("J104").Value = 25112
$ws.Range("L104").Value = 25112
$ws.Range("N104").Value = -32100
$ws.Range("H116").Value = 2840.2563
$ws.Range("I116").Value = 1479.1364
$ws.Range("K116").Value = 1479.1364
$ws.Range("M116").Value = 814.8635999999999
$ws.Range("H125").Value = 77743.42999999999
$ws.Range("J125").Value = 77743.42999999999
$ws.Range("L125").Value = 77743.42999999999
$ws.Range("N125").Value = -87583.42999999999
$ws.Range("H132").Value = 3035450.8
$ws.Range("I132").Value = 3338473
$ws.Range("J132").Value = 5229.6665
$ws.Range("K132").Value = 10015419
$ws.Range("L132").Value = 15688.9995
$ws.Range("M132").Value = -10012889
$ws.Range("N132").Value = -20748.9995
$ws.Range("H136").Value = 50457950
$ws.Range("I136").Value = 50457950
$ws.Range("K136").Value = 151373850
$ws.Range("M136").Value = -151371300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2840.2563
$ws.Range("I3").Value = 1479.1364
$ws.Range("K3").Value = 1479.1364
$ws.Range("M3").Value = -1365.1364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7213.0225
$ws.Range("I31").Value = 2278.5
$ws.Range("K31").Value = 2278.5
$ws.Range("M31").Value = -1983.5
$ws.Range("H34").Value = 7213.0225
$ws.Range("I34").Value = 2278.5
$ws.Range("K34").Value = 2278.5
$ws.Range("M34").Value = -2076.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 214.28572
$ws.Range("I2").Value = 191.22223
$ws.Range("K2").Value = 1147.33338
$ws.Range("M2").Value = -1034.33338
$ws.Range("H11").Value = 921496.4399999999
$ws.Range("I11").Value = 959850.4399999999
$ws.Range("K11").Value = 2879551.32
$ws.Range("M11").Value = -2879411.32
$ws.Range("H87").Value = 7999
$ws.Range("I87").Value = 7999
$ws.Range("K87").Value = 23997
$ws.Range("M87").Value = -22749
$ws.Range("H90").Value = 7999
$ws.Range("I90").Value = 7999
$ws.Range("K90").Value = 71991
$ws.Range("M90").Value = -65751
$ws.Range("H107").Value = 1455.7407
$ws.Range("J107").Value = 1987.6666
$ws.Range("L107").Value = 5962.9998
$ws.Range("N107").Value = -9802.9998
$ws.Range("H120").Value = 18990
$ws.Range("I120").Value = 18990
$ws.Range("K120").Value = 56970
$ws.Range("M120").Value = -52132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1527.2142
$ws.Range("I97").Value = 1358.3636
$ws.Range("J97").Value = 2146.3333
$ws.Range("K97").Value = 1358.3636
$ws.Range("L97").Value = 2146.3333
$ws.Range("M97").Value = -862.3635999999999
$ws.Range("N97").Value = -3138.3333
$ws.Range("H109").Value = 64950
$ws.Range("J109").Value = 64950
$ws.Range("L109").Value = 64950
$ws.Range("N109").Value = -67030
$ws.Range("H126").Value = 6741.25
$ws.Range("I126").Value = 6963.636
$ws.Range("J126").Value = 4295
$ws.Range("K126").Value = 20890.908
$ws.Range("L126").Value = 12885
$ws.Range("M126").Value = -18420.908
$ws.Range("N126").Value = -17825
$ws.Range("H132").Value = 7814528.5
$ws.Range("I132").Value = 9616849
$ws.Range("J132").Value = 4471.3335
$ws.Range("K132").Value = 28850547
$ws.Range("L132").Value = 13414.0005
$ws.Range("M132").Value = -28848017
$ws.Range("N132").Value = -18474.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1589.92
$ws.Range("J22").Value = 2274.8333
$ws.Range("L22").Value = 2274.8333
$ws.Range("N22").Value = -2864.8333
$ws.Range("H27").Value = 1589.92
$ws.Range("J27").Value = 2274.8333
$ws.Range("L27").Value = 2274.8333
$ws.Range("N27").Value = -2488.8333
$ws.Range("H31").Value = 10059.5
$ws.Range("I31").Value = 9049.5
$ws.Range("J31").Value = 11574.5
$ws.Range("K31").Value = 9049.5
$ws.Range("L31").Value = 11574.5
$ws.Range("M31").Value = -8801.5
$ws.Range("N31").Value = -12070.5
$ws.Range("H43").Value = 49991.332
$ws.Range("J43").Value = 49987
$ws.Range("L43").Value = 49987
$ws.Range("N43").Value = -50373
$ws.Range("H46").Value = 807.6667
$ws.Range("J46").Value = 933.6667
$ws.Range("L46").Value = 933.6667
$ws.Range("N46").Value = -1309.6667
$ws.Range("H68").Value = 126118.5
$ws.Range("I68").Value = 1278.4286
$ws.Range("J68").Value = 999999
$ws.Range("K68").Value = 1278.4286
$ws.Range("L68").Value = 999999
$ws.Range("M68").Value = -529.4286
$ws.Range("N68").Value = -1001497
$ws.Range("H71").Value = 126118.5
$ws.Range("I71").Value = 1278.4286
$ws.Range("J71").Value = 999999
$ws.Range("K71").Value = 6392.143
$ws.Range("L71").Value = 4999995
$ws.Range("M71").Value = -2648.143
$ws.Range("N71").Value = -5007483

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1098.1428
$ws.Range("I100").Value = 1062.2273
$ws.Range("K100").Value = 2124.4546
$ws.Range("M100").Value = -1583.4546
$ws.Range("H122").Value = 6422.25
$ws.Range("I122").Value = 6986.8
$ws.Range("J122").Value = 3599.5
$ws.Range("K122").Value = 20960.4
$ws.Range("L122").Value = 10798.5
$ws.Range("M122").Value = -18510.4
$ws.Range("N122").Value = -15698.5
$ws.Range("H126").Value = 1829.2941
$ws.Range("I126").Value = 1857.2142
$ws.Range("J126").Value = 1699
$ws.Range("K126").Value = 5571.642599999999
$ws.Range("L126").Value = 5097
$ws.Range("M126").Value = -3101.642599999999
$ws.Range("N126").Value = -10037
$ws.Range("H132").Value = 11911347
$ws.Range("I132").Value = 15156335
$ws.Range("K132").Value = 45469005
$ws.Range("M132").Value = -45466475
$ws.Range("H136").Value = 20001356
$ws.Range("I136").Value = 22728020
$ws.Range("J136").Value = 5816.6665
$ws.Range("K136").Value = 68184060
$ws.Range("L136").Value = 17449.9995
$ws.Range("M136").Value = -68181510
$ws.Range("N136").Value = -22549.9995
$ws.Range("H139").Value = 144374.38
$ws.Range("J139").Value = 144374.38
$ws.Range("L139").Value = 144374.38
$ws.Range("N139").Value = -154654.38

Write-Output "Applied all updates"